# Rebuild Sheet1 header row (No/Type/Pertanyaan/Aspek/Kriteria) and add a
# new Sheet2 with an "aspek/kriteria/bobot" 1-5 scoring grid.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"

# --- Sheet1: update header row text ---
$ws1.Range("A1").Value = "No"
$ws1.Range("B1").Value = "Type"
$ws1.Range("C1").Value = "Pertanyaan"
$ws1.Range("D1").Value = "Aspek"
$ws1.Range("E1").Value = "Kriteria"

$ws1.Range("A1:E1").HorizontalAlignment = -4108  # xlCenter

$ws1.Columns.Item(2).ColumnWidth = 13.5
$ws1.Columns.Item(3).ColumnWidth = 34.333333333333336
$ws1.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws1.Columns.Item(5).ColumnWidth = 16.5

# --- Sheet2: clone Sheet1 (keeps the same sheetFormatPr/view defaults)
#     then overwrite its content entirely ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

$ws2.Cells.ClearContents()
$ws2.Cells.ClearFormats()

$ws2.Range("A1").Value = "no"
$ws2.Range("B1").Value = "aspek"
$ws2.Range("C1").Value = "kriteria"
$ws2.Range("D1").Value = "bobot"

$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 2
$ws2.Range("F2").Value = 3
$ws2.Range("G2").Value = 4
$ws2.Range("H2").Value = 5

$ws2.Range("A1:H2").HorizontalAlignment = -4108  # xlCenter

$ws2.Range("D1:H1").Merge()
$ws2.Range("A1:A2").Merge()
$ws2.Range("B1:B2").Merge()
$ws2.Range("C1:C2").Merge()

$ws2.Columns.Item(2).ColumnWidth = 13.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(4).ColumnWidth = 25.666666666666668
$ws2.Columns.Item(5).ColumnWidth = 19.833333333333332
$ws2.Columns.Item(6).ColumnWidth = 19.333333333333332
$ws2.Columns.Item(7).ColumnWidth = 19.333333333333332
$ws2.Columns.Item(8).ColumnWidth = 23.5

# --- view state: Sheet1 zoomed to 140%, selection A2:E3 ---
$ws1.Activate()
$excel.ActiveWindow.Zoom = 140
$ws1.Range("A2:E3").Select()

# --- view state: Sheet2 active/zoomed to 90%, selection H12 ---
$ws2.Activate()
$excel.ActiveWindow.Zoom = 90
$ws2.Range("H12").Select()
